$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 333335400
$ws.Range("I76").Value = 333335400
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 333335400
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -333335085
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 333335400
$ws.Range("I79").Value = 333335400
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 333335400
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -333334308
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 2526251
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 3367668
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 3367668
$ws.Range("M86").Value = -877
$ws.Range("N86").Value = -3369914
$ws.Range("H88").Value = 1870.2858
$ws.Range("I88").Value = 993
$ws.Range("J88").Value = 2221.2
$ws.Range("K88").Value = 993
$ws.Range("L88").Value = 2221.2
$ws.Range("M88").Value = -587
$ws.Range("N88").Value = -3033.2
$ws.Range("H89").Value = 2526251
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 3367668
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 16838340
$ws.Range("M89").Value = -4384
$ws.Range("N89").Value = -16849572
$ws.Range("H91").Value = 1870.2858
$ws.Range("I91").Value = 993
$ws.Range("J91").Value = 2221.2
$ws.Range("K91").Value = 993
$ws.Range("L91").Value = 2221.2
$ws.Range("M91").Value = 411
$ws.Range("N91").Value = -5029.2
$ws.Range("H118").Value = 1823.421
$ws.Range("I118").Value = 567.5
$ws.Range("J118").Value = 2736.818
$ws.Range("K118").Value = 1702.5
$ws.Range("L118").Value = 8210.454000000002
$ws.Range("M118").Value = -45.5
$ws.Range("N118").Value = -11524.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1741703.8
$ws.Range("I32").Value = 1989296.9
$ws.Range("K32").Value = 1989296.9
$ws.Range("M32").Value = -1989009.9
$ws.Range("H61").Value = 5912200
$ws.Range("I61").Value = 2778738.2
$ws.Range("J61").Value = 29413164
$ws.Range("K61").Value = 2778738.2
$ws.Range("L61").Value = 29413164
$ws.Range("M61").Value = -2778526.2
$ws.Range("N61").Value = -29413588
$ws.Range("H132").Value = 13893795
$ws.Range("J132").Value = 9262426
$ws.Range("L132").Value = 27787278
$ws.Range("N132").Value = -27792338
$ws.Range("H136").Value = 5912200
$ws.Range("I136").Value = 2778738.2
$ws.Range("J136").Value = 29413164
$ws.Range("K136").Value = 8336214.600000001
$ws.Range("L136").Value = 88239492
$ws.Range("M136").Value = -8333664.600000001
$ws.Range("N136").Value = -88244592

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 14479595
$ws.Range("I134").Value = 14706737
$ws.Range("J134").Value = 11905328
$ws.Range("K134").Value = 44120211
$ws.Range("L134").Value = 35715984
$ws.Range("M134").Value = -44117676
$ws.Range("N134").Value = -35721054

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2659946.5
$ws.Range("I5").Value = 2263119.2
$ws.Range("J5").Value = 3334553
$ws.Range("K5").Value = 6789357.600000001
$ws.Range("L5").Value = 10003659
$ws.Range("M5").Value = -6789245.600000001
$ws.Range("N5").Value = -10003883
$ws.Range("H122").Value = 590.9167
$ws.Range("J122").Value = 1899.5
$ws.Range("L122").Value = 17095.5
$ws.Range("N122").Value = -21995.5
$ws.Range("H135").Value = 2659946.5
$ws.Range("I135").Value = 2263119.2
$ws.Range("J135").Value = 3334553
$ws.Range("K135").Value = 20368072.8
$ws.Range("L135").Value = 30010977
$ws.Range("M135").Value = -20365537.8
$ws.Range("N135").Value = -30016047
$ws.Range("H141").Value = 2484.6155
$ws.Range("I141").Value = 2484.6155
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7453.8465
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2273.8465
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13514.571
$ws.Range("I80").Value = 5541.6665
$ws.Range("J80").Value = 24145.111
$ws.Range("K80").Value = 5541.6665
$ws.Range("L80").Value = 24145.111
$ws.Range("M80").Value = -4543.6665
$ws.Range("N80").Value = -26141.111
$ws.Range("H83").Value = 13514.571
$ws.Range("I83").Value = 5541.6665
$ws.Range("J83").Value = 24145.111
$ws.Range("K83").Value = 27708.3325
$ws.Range("L83").Value = 120725.555
$ws.Range("M83").Value = -22716.3325
$ws.Range("N83").Value = -130709.555
$ws.Range("H113").Value = 25309.908
$ws.Range("I113").Value = 3544.1428
$ws.Range("K113").Value = 3544.1428
$ws.Range("M113").Value = -1374.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1521.3572
$ws.Range("I61").Value = 1124.75
$ws.Range("J61").Value = 1680
$ws.Range("K61").Value = 1124.75
$ws.Range("L61").Value = 1680
$ws.Range("M61").Value = -922.75
$ws.Range("N61").Value = -2084
$ws.Range("H68").Value = 3373.5
$ws.Range("I68").Value = 3306.923
$ws.Range("J68").Value = 3497.1428
$ws.Range("K68").Value = 3306.923
$ws.Range("L68").Value = 3497.1428
$ws.Range("M68").Value = -2557.923
$ws.Range("N68").Value = -4995.1428
$ws.Range("H71").Value = 3373.5
$ws.Range("I71").Value = 3306.923
$ws.Range("J71").Value = 3497.1428
$ws.Range("K71").Value = 16534.615
$ws.Range("L71").Value = 17485.714
$ws.Range("M71").Value = -12790.615
$ws.Range("N71").Value = -24973.714
$ws.Range("H113").Value = 1521.3572
$ws.Range("I113").Value = 1124.75
$ws.Range("J113").Value = 1680
$ws.Range("K113").Value = 1124.75
$ws.Range("L113").Value = 1680
$ws.Range("M113").Value = 1045.25
$ws.Range("N113").Value = -6020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H81").Value = 12481.556
$ws.Range("I81").Value = 1429
$ws.Range("J81").Value = 16349.95
$ws.Range("K81").Value = 2858
$ws.Range("L81").Value = 32699.9
$ws.Range("M81").Value = -1797
$ws.Range("N81").Value = -34821.9
$ws.Range("H84").Value = 12481.556
$ws.Range("I84").Value = 1429
$ws.Range("J84").Value = 16349.95
$ws.Range("K84").Value = 14290
$ws.Range("L84").Value = 163499.5
$ws.Range("M84").Value = -8986
$ws.Range("N84").Value = -174107.5
